$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.419591
$ws.Range("H2").Value = 4.258773
$ws.Range("I2").Value = 0.001848767113890483
$ws.Range("J2").Value = 0.001848767113890483
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.229822
$ws.Range("N2").Value = 0.689466
$ws.Range("O2").Value = 0.09226175421862418
$ws.Range("P2").Value = 0.09226175421862419
$ws.Range("Q2").Value = 0.326253242802
$ws.Range("R2").Value = 2.936279185218
$ws.Range("S2").Value = 0.0001705704970692389
$ws.Range("T2").Value = 0.000170570497069239
$ws.Range("G3").Value = 1.419591
$ws.Range("H3").Value = 4.258773
$ws.Range("I3").Value = 0.001848767113890483
$ws.Range("J3").Value = 0.001848767113890483
$ws.Range("O3").Value = 0.4364142651333466
$ws.Range("P3").Value = 0.4364142651333466
$ws.Range("Q3").Value = 1.543235010115
$ws.Range("R3").Value = 13.889115091035
$ws.Range("S3").Value = 0.0008068283414112133
$ws.Range("T3").Value = 0.0008068283414112133
$ws.Range("G4").Value = 1.419591
$ws.Range("H4").Value = 4.258773
$ws.Range("I4").Value = 0.001848767113890483
$ws.Range("J4").Value = 0.001848767113890483
$ws.Range("M4").Value = 1.174057666666666
$ws.Range("O4").Value = 0.4713239806480292
$ws.Range("P4").Value = 0.4713239806480293
$ws.Range("Q4").Value = 1.666681697080999
$ws.Range("R4").Value = 15.000135273729
$ws.Range("S4").Value = 0.000871368275410031
$ws.Range("T4").Value = 0.0008713682754100311
$ws.Range("I5").Value = 0.9578582377148513
$ws.Range("J5").Value = 0.9578582377148513
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.229822
$ws.Range("N5").Value = 0.689466
$ws.Range("O5").Value = 0.09226175421862418
$ws.Range("P5").Value = 0.09226175421862419
$ws.Range("Q5").Value = 169.0339220397833
$ws.Range("R5").Value = 1521.30529835805
$ws.Range("S5").Value = 0.08837368130433211
$ws.Range("T5").Value = 0.08837368130433212
$ws.Range("I6").Value = 0.9578582377148513
$ws.Range("J6").Value = 0.9578582377148513
$ws.Range("O6").Value = 0.4364142651333466
$ws.Range("P6").Value = 0.4364142651333466
$ws.Range("Q6").Value = 799.5600722569861
$ws.Range("R6").Value = 7196.040650312874
$ws.Range("S6").Value = 0.4180229989142493
$ws.Range("T6").Value = 0.4180229989142493
$ws.Range("I7").Value = 0.9578582377148513
$ws.Range("J7").Value = 0.9578582377148513
$ws.Range("M7").Value = 1.174057666666666
$ws.Range("O7").Value = 0.4713239806480292
$ws.Range("P7").Value = 0.4713239806480293
$ws.Range("Q7").Value = 863.5186017767804
$ws.Range("R7").Value = 7771.667415991024
$ws.Range("S7").Value = 0.45146155749627
$ws.Range("T7").Value = 0.4514615574962701
$ws.Range("H8").Value = 92.81792100000001
$ws.Range("I8").Value = 0.04029299517125823
$ws.Range("J8").Value = 0.04029299517125823
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.229822
$ws.Range("N8").Value = 0.689466
$ws.Range("O8").Value = 0.09226175421862418
$ws.Range("P8").Value = 0.09226175421862419
$ws.Range("Q8").Value = 7.110533413354001
$ws.Range("R8").Value = 63.99480072018601
$ws.Range("S8").Value = 0.003717502417222837
$ws.Range("T8").Value = 0.003717502417222838
$ws.Range("H9").Value = 92.81792100000001
$ws.Range("I9").Value = 0.04029299517125823
$ws.Range("J9").Value = 0.04029299517125823
$ws.Range("O9").Value = 0.4364142651333466
$ws.Range("P9").Value = 0.4364142651333466
$ws.Range("Q9").Value = 33.63406907418833
$ws.Range("R9").Value = 302.706621667695
$ws.Range("S9").Value = 0.01758443787768614
$ws.Range("T9").Value = 0.01758443787768614
$ws.Range("H10").Value = 92.81792100000001
$ws.Range("I10").Value = 0.04029299517125823
$ws.Range("J10").Value = 0.04029299517125823
$ws.Range("M10").Value = 1.174057666666666
$ws.Range("O10").Value = 0.4713239806480292
$ws.Range("P10").Value = 0.4713239806480293
$ws.Range("S10").Value = 0.01899105487634925
$ws.Range("T10").Value = 0.01899105487634925

Write-Host "Applied TPM updates"
